# Auto-generated Excel COM-interop script applying numeric updates
# to the Anima Profits tracker sheets (per-sheet "Sheets/Anima_Profits.xlsx" diff).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 998.2
$ws.Range("I32").Value = 950.25
$ws.Range("J32").Value = 1015.63635
$ws.Range("K32").Value = 950.25
$ws.Range("L32").Value = 1015.63635
$ws.Range("M32").Value = -624.25
$ws.Range("N32").Value = -1667.63635

$ws.Range("H53").Value = 226.46153
$ws.Range("I53").Value = 154.2
$ws.Range("J53").Value = 467.33334
$ws.Range("K53").Value = 154.2
$ws.Range("L53").Value = 467.33334
$ws.Range("M53").Value = 482.8
$ws.Range("N53").Value = -1741.33334

$ws.Range("H57").Value = 51000
$ws.Range("J57").Value = 51000
$ws.Range("L57").Value = 153000
$ws.Range("N57").Value = -153998

$ws.Range("H138").Value = 1947.8918
$ws.Range("J138").Value = 2725.111
$ws.Range("L138").Value = 8175.333
$ws.Range("N138").Value = -18455.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2750.5833
$ws.Range("I61").Value = 1734.3334
$ws.Range("J61").Value = 4444.3335
$ws.Range("K61").Value = 1734.3334
$ws.Range("L61").Value = 4444.3335
$ws.Range("M61").Value = -1522.3334
$ws.Range("N61").Value = -4868.3335

$ws.Range("H74").Value = 1400.5333
$ws.Range("I74").Value = 997.2727
$ws.Range("J74").Value = 2509.5
$ws.Range("K74").Value = 997.2727
$ws.Range("L74").Value = 2509.5
$ws.Range("M74").Value = -123.2727
$ws.Range("N74").Value = -4257.5

$ws.Range("H77").Value = 1400.5333
$ws.Range("I77").Value = 997.2727
$ws.Range("J77").Value = 2509.5
$ws.Range("K77").Value = 4986.363499999999
$ws.Range("L77").Value = 12547.5
$ws.Range("M77").Value = -618.3634999999995
$ws.Range("N77").Value = -21283.5

$ws.Range("H103").Value = 95000
$ws.Range("J103").Value = 95000
$ws.Range("L103").Value = 95000
$ws.Range("N103").Value = -97344

$ws.Range("H122").Value = 2043.9445
$ws.Range("I122").Value = 1922.2307
$ws.Range("J122").Value = 2360.4
$ws.Range("K122").Value = 5766.6921
$ws.Range("L122").Value = 7081.200000000001
$ws.Range("M122").Value = -3316.6921
$ws.Range("N122").Value = -11981.2

$ws.Range("H136").Value = 2750.5833
$ws.Range("I136").Value = 1734.3334
$ws.Range("J136").Value = 4444.3335
$ws.Range("K136").Value = 5203.0002
$ws.Range("L136").Value = 13333.0005
$ws.Range("M136").Value = -2653.0002
$ws.Range("N136").Value = -18433.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 19000
$ws.Range("J74").Value = 19000
$ws.Range("L74").Value = 19000
$ws.Range("N74").Value = -20872

$ws.Range("H77").Value = 19000
$ws.Range("J77").Value = 19000
$ws.Range("L77").Value = 57000
$ws.Range("N77").Value = -66360

$ws.Range("H86").Value = 200002400
$ws.Range("I86").Value = 250002240
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 250002240
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -250001117
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 200002400
$ws.Range("I89").Value = 250002240
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 1250011200
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -1250005584
$ws.Range("N89").Value = -26232

$ws.Range("H134").Value = 3991.3157
$ws.Range("I134").Value = 4445.364
$ws.Range("K134").Value = 13336.092
$ws.Range("M134").Value = -10801.092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 29889.111
$ws.Range("J4").Value = 29889.111
$ws.Range("L4").Value = 29889.111
$ws.Range("N4").Value = -30113.111

$ws.Range("H31").Value = 7111.978
$ws.Range("I31").Value = 1393.6316
$ws.Range("J31").Value = 11136
$ws.Range("K31").Value = 1393.6316
$ws.Range("L31").Value = 11136
$ws.Range("M31").Value = -1098.6316
$ws.Range("N31").Value = -11726

$ws.Range("H34").Value = 7111.978
$ws.Range("I34").Value = 1393.6316
$ws.Range("J34").Value = 11136
$ws.Range("K34").Value = 1393.6316
$ws.Range("L34").Value = 11136
$ws.Range("M34").Value = -1191.6316
$ws.Range("N34").Value = -11540

$ws.Range("H134").Value = 2413.3914
$ws.Range("I134").Value = 1452.7894
$ws.Range("K134").Value = 4358.3682
$ws.Range("M134").Value = -1823.3682

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""

$ws.Range("H6").Value = 295.35715
$ws.Range("I6").Value = 41.75
$ws.Range("J6").Value = 633.5
$ws.Range("K6").Value = 125.25
$ws.Range("L6").Value = 1900.5
$ws.Range("M6").Value = -12.25
$ws.Range("N6").Value = -2126.5

$ws.Range("H38").Value = 1566
$ws.Range("I38").Value = 71.59999999999999
$ws.Range("J38").Value = 2500
$ws.Range("K38").Value = 214.8
$ws.Range("L38").Value = 7500
$ws.Range("M38").Value = 132.2
$ws.Range("N38").Value = -8194

$ws.Range("H86").Value = 1140
$ws.Range("J86").Value = 1980
$ws.Range("L86").Value = 5940
$ws.Range("N86").Value = -8312

$ws.Range("H89").Value = 1140
$ws.Range("J89").Value = 1980
$ws.Range("L89").Value = 17820
$ws.Range("N89").Value = -29676

$ws.Range("H99").Value = 200
$ws.Range("I99").Value = 200
$ws.Range("K99").Value = 600
$ws.Range("M99").Value = 1646

$ws.Range("H108").Value = 1142
$ws.Range("I108").Value = 770.4
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 2311.2
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 568.8000000000002
$ws.Range("N108").Value = -14760

$ws.Range("H120").Value = 9036.666999999999
$ws.Range("I120").Value = 3110
$ws.Range("K120").Value = 9330
$ws.Range("M120").Value = -4492

$ws.Range("H132").Value = 1943.5555
$ws.Range("I132").Value = 1626.6666
$ws.Range("J132").Value = 2102
$ws.Range("K132").Value = 14639.9994
$ws.Range("L132").Value = 18918
$ws.Range("M132").Value = -12109.9994
$ws.Range("N132").Value = -23978

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 34000
$ws.Range("J15").Value = 34000
$ws.Range("L15").Value = 34000
$ws.Range("N15").Value = -34576

$ws.Range("H81").Value = 34000
$ws.Range("J81").Value = 34000
$ws.Range("L81").Value = 34000
$ws.Range("N81").Value = -35996

$ws.Range("H84").Value = 34000
$ws.Range("J84").Value = 34000
$ws.Range("L84").Value = 102000
$ws.Range("N84").Value = -111984

$ws.Range("H102").Value = 1609.75
$ws.Range("I102").Value = 1528.2142
$ws.Range("J102").Value = 1800
$ws.Range("K102").Value = 1528.2142
$ws.Range("L102").Value = 1800
$ws.Range("M102").Value = 93.78580000000011
$ws.Range("N102").Value = -5044

$ws.Range("H105").Value = 33000
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 33000
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 33000
$ws.Range("M105").Value = ""
$ws.Range("N105").Value = -39988

$ws.Range("H126").Value = 1719
$ws.Range("I126").Value = 1722.4546
$ws.Range("J126").Value = 1700
$ws.Range("K126").Value = 5167.3638
$ws.Range("L126").Value = 5100
$ws.Range("M126").Value = -2697.3638
$ws.Range("N126").Value = -10040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2665
$ws.Range("I7").Value = 2275
$ws.Range("J7").Value = 5005
$ws.Range("K7").Value = 2275
$ws.Range("L7").Value = 5005
$ws.Range("M7").Value = -2163
$ws.Range("N7").Value = -5229

$ws.Range("H126").Value = 2665
$ws.Range("I126").Value = 2275
$ws.Range("J126").Value = 5005
$ws.Range("K126").Value = 6825
$ws.Range("L126").Value = 15015
$ws.Range("M126").Value = -4355
$ws.Range("N126").Value = -19955

$ws.Range("H136").Value = 9805934
$ws.Range("I136").Value = 2268.6667
$ws.Range("J136").Value = 33334730
$ws.Range("K136").Value = 6806.000100000001
$ws.Range("L136").Value = 100004190
$ws.Range("M136").Value = -4256.000100000001
$ws.Range("N136").Value = -100009290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 39750.75
$ws.Range("J4").Value = 39750.75
$ws.Range("L4").Value = 39750.75
$ws.Range("N4").Value = -39976.75

$ws.Range("H100").Value = 1657.3334
$ws.Range("I100").Value = 1963.125
$ws.Range("J100").Value = 1045.75
$ws.Range("K100").Value = 3926.25
$ws.Range("L100").Value = 2091.5
$ws.Range("M100").Value = -3385.25
$ws.Range("N100").Value = -3173.5

$ws.Range("H122").Value = 2413.5386
$ws.Range("I122").Value = 2013.909
$ws.Range("J122").Value = 2706.6
$ws.Range("K122").Value = 6041.727000000001
$ws.Range("L122").Value = 8119.799999999999
$ws.Range("M122").Value = -3591.727000000001
$ws.Range("N122").Value = -13019.8

$ws.Range("H136").Value = 2146.2354
$ws.Range("I136").Value = 1874.8055
$ws.Range("K136").Value = 5624.416499999999
$ws.Range("M136").Value = -3074.416499999999

